# "Added quiz game creator"
#
# This edit mostly reshuffles view-state (selections / active tab) across
# several sheets and fills in a couple of previously-empty cells that back
# the new quiz game feature (BirdSoundsURL!B6:C6, and a re-pointed
# BirdTag!B8). The shared string "medium" -> "TEMP" rename falls naturally
# out of re-pointing BirdTag!B8 away from it and writing BirdSoundsURL!B6
# with the final "TEMP" text, since that cell ends up the sole remaining
# reference to that shared-string slot.

$wb = $excel.ActiveWorkbook

# --- BirdTag: point B8 at "small" instead of "medium" -----------------
$wsTag = $wb.Worksheets.Item("BirdTag")
$wsTag.Range("B8").Value = "small"

# --- BirdSoundsURL: fill in the previously-empty row 6 (aegithalos
#     caudatus) and move the selection ----------------------------------
$wsSounds = $wb.Worksheets.Item("BirdSoundsURL")
$wsSounds.Range("B6").Value = "TEMP"
$wsSounds.Range("C6").Value = "mixed"
$wsSounds.Range("F8").Select()

# --- Language: move selection -------------------------------------------
$wsLang = $wb.Worksheets.Item("Language")
$wsLang.Range("H10").Select()

# --- BirdOrderTranslation: widen column A and move selection -----------
$wsOrderTr = $wb.Worksheets.Item("BirdOrderTranslation")
$wsOrderTr.Columns("A").ColumnWidth = 14.29
$wsOrderTr.Range("G11").Select()

# --- BirdNameTranslation: move selection --------------------------------
$wsNameTr = $wb.Worksheets.Item("BirdNameTranslation")
$wsNameTr.Range("C21").Select()

# --- BirdTag: final active sheet/tab + selection ------------------------
# Selected last so it ends up as the workbook's active tab/sheet.
$wsTag.Range("B8").Select()

Write-Output "done"
